# Generate Report for Handback
# Update the timestamp cells that record when the handoff/handback xliff
# files were generated, reflecting a newer report generation run.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" for the first file (shared text between
# the Overview sheet and the de-de detail sheet).
$wsOverview.Range("G2").Value = "2016-09-04 01:10:46"
$wsDeDe.Range("H2").Value = "2016-09-04 01:10:46"

# zh-cn detail sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-04 01:10:35"
$wsZhCn.Range("K2").Value = "2016-09-04 01:11:08"

# de-de detail sheet: Correspond Handback DateTime
$wsDeDe.Range("K2").Value = "2016-09-04 01:11:16"
